$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Cells.Item(18, 8).Value = 1362.1052
$ws.Cells.Item(18, 9).Value = 1265.2
$ws.Cells.Item(18, 10).Value = 1725.5
$ws.Cells.Item(18, 11).Value = 1265.2
$ws.Cells.Item(18, 12).Value = 1725.5
$ws.Cells.Item(18, 13).Value = -981.2
$ws.Cells.Item(18, 14).Value = -2293.5
# Row 75
$ws.Cells.Item(75, 8).Value = 26527.143
$ws.Cells.Item(75, 10).Value = 28615
$ws.Cells.Item(75, 12).Value = 28615
$ws.Cells.Item(75, 14).Value = -30487
# Row 78
$ws.Cells.Item(78, 8).Value = 26527.143
$ws.Cells.Item(78, 10).Value = 28615
$ws.Cells.Item(78, 12).Value = 85845
$ws.Cells.Item(78, 14).Value = -95205
# Row 100
$ws.Cells.Item(100, 8).Value = 17780066
$ws.Cells.Item(100, 9).Value = 29631870
$ws.Cells.Item(100, 10).Value = 2360
$ws.Cells.Item(100, 11).Value = 29631870
$ws.Cells.Item(100, 12).Value = 2360
$ws.Cells.Item(100, 13).Value = -29631329
$ws.Cells.Item(100, 14).Value = -3442
# Row 132
$ws.Cells.Item(132, 8).Value = 4187154.8
$ws.Cells.Item(132, 9).Value = 5581493.5
$ws.Cells.Item(132, 10).Value = 4138
$ws.Cells.Item(132, 11).Value = 16744480.5
$ws.Cells.Item(132, 12).Value = 12414
$ws.Cells.Item(132, 13).Value = -16741950.5
$ws.Cells.Item(132, 14).Value = -17474
# Row 137
$ws.Cells.Item(137, 8).Value = 1247.1875
$ws.Cells.Item(137, 9).Value = 781.8148
$ws.Cells.Item(137, 11).Value = 2345.4444
$ws.Cells.Item(137, 13).Value = 204.5556000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 6413244.5
$ws.Cells.Item(2, 9).Value = 8336847
$ws.Cells.Item(2, 10).Value = 1237.6666
$ws.Cells.Item(2, 11).Value = 8336847
$ws.Cells.Item(2, 12).Value = 1237.6666
$ws.Cells.Item(2, 13).Value = -8336734
$ws.Cells.Item(2, 14).Value = -1463.6666
# Row 32
$ws.Cells.Item(32, 8).Value = 16510.467
$ws.Cells.Item(32, 9).Value = 12607.782
$ws.Cells.Item(32, 10).Value = 59440
$ws.Cells.Item(32, 11).Value = 12607.782
$ws.Cells.Item(32, 12).Value = 59440
$ws.Cells.Item(32, 13).Value = -12320.782
$ws.Cells.Item(32, 14).Value = -60014
# Row 33
$ws.Cells.Item(33, 8).Value = 461713.8
$ws.Cells.Item(33, 9).Value = 1252388
$ws.Cells.Item(33, 10).Value = 9900
$ws.Cells.Item(33, 11).Value = 1252388
$ws.Cells.Item(33, 12).Value = 9900
$ws.Cells.Item(33, 13).Value = -1252059
$ws.Cells.Item(33, 14).Value = -10558
# Row 36
$ws.Cells.Item(36, 8).Value = 4692.467
$ws.Cells.Item(36, 9).Value = 898
$ws.Cells.Item(36, 10).Value = 8012.625
$ws.Cells.Item(36, 11).Value = 898
$ws.Cells.Item(36, 12).Value = 8012.625
$ws.Cells.Item(36, 13).Value = -552
$ws.Cells.Item(36, 14).Value = -8704.625
# Row 61
$ws.Cells.Item(61, 8).Value = 2520.4856
$ws.Cells.Item(61, 9).Value = 2171.5417
$ws.Cells.Item(61, 10).Value = 3281.818
$ws.Cells.Item(61, 11).Value = 2171.5417
$ws.Cells.Item(61, 12).Value = 3281.818
$ws.Cells.Item(61, 13).Value = -1959.5417
$ws.Cells.Item(61, 14).Value = -3705.818
# Row 98
$ws.Cells.Item(98, 8).Value = 32490
$ws.Cells.Item(98, 10).Value = 32490
$ws.Cells.Item(98, 12).Value = 32490
$ws.Cells.Item(98, 14).Value = -38480
# Row 105
$ws.Cells.Item(105, 8).Value = 29870
$ws.Cells.Item(105, 10).Value = 29870
$ws.Cells.Item(105, 12).Value = 29870
$ws.Cells.Item(105, 14).Value = -36858
# Row 116
$ws.Cells.Item(116, 8).Value = 6413244.5
$ws.Cells.Item(116, 9).Value = 8336847
$ws.Cells.Item(116, 10).Value = 1237.6666
$ws.Cells.Item(116, 11).Value = 8336847
$ws.Cells.Item(116, 12).Value = 1237.6666
$ws.Cells.Item(116, 13).Value = -8334553
$ws.Cells.Item(116, 14).Value = -5825.6666
# Row 122
$ws.Cells.Item(122, 8).Value = 6945800.5
$ws.Cells.Item(122, 9).Value = 9616389
$ws.Cells.Item(122, 10).Value = 2270
$ws.Cells.Item(122, 11).Value = 28849167
$ws.Cells.Item(122, 12).Value = 6810
$ws.Cells.Item(122, 13).Value = -28846717
$ws.Cells.Item(122, 14).Value = -11710
# Row 136
$ws.Cells.Item(136, 8).Value = 2520.4856
$ws.Cells.Item(136, 9).Value = 2171.5417
$ws.Cells.Item(136, 10).Value = 3281.818
$ws.Cells.Item(136, 11).Value = 6514.625100000001
$ws.Cells.Item(136, 12).Value = 9845.454000000002
$ws.Cells.Item(136, 13).Value = -3964.625100000001
$ws.Cells.Item(136, 14).Value = -14945.454

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 6413244.5
$ws.Cells.Item(3, 9).Value = 8336847
$ws.Cells.Item(3, 10).Value = 1237.6666
$ws.Cells.Item(3, 11).Value = 8336847
$ws.Cells.Item(3, 12).Value = 1237.6666
$ws.Cells.Item(3, 13).Value = -8336733
$ws.Cells.Item(3, 14).Value = -1465.6666
# Row 55
$ws.Cells.Item(55, 8).Value = 62526.668
$ws.Cells.Item(55, 10).Value = 62526.668
$ws.Cells.Item(55, 12).Value = 62526.668
$ws.Cells.Item(55, 14).Value = -63072.668
# Row 75
$ws.Cells.Item(75, 8).Value = 13808.363
$ws.Cells.Item(75, 9).Value = 8788
$ws.Cells.Item(75, 10).Value = 36400
$ws.Cells.Item(75, 11).Value = 8788
$ws.Cells.Item(75, 12).Value = 36400
$ws.Cells.Item(75, 13).Value = -7852
$ws.Cells.Item(75, 14).Value = -38272
# Row 78
$ws.Cells.Item(78, 8).Value = 13808.363
$ws.Cells.Item(78, 9).Value = 8788
$ws.Cells.Item(78, 10).Value = 36400
$ws.Cells.Item(78, 11).Value = 26364
$ws.Cells.Item(78, 12).Value = 109200
$ws.Cells.Item(78, 13).Value = -21684
$ws.Cells.Item(78, 14).Value = -118560

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 60351.89
$ws.Cells.Item(16, 9).Value = 169670.33
$ws.Cells.Item(16, 10).Value = 5692.6665
$ws.Cells.Item(16, 11).Value = 169670.33
$ws.Cells.Item(16, 12).Value = 5692.6665
$ws.Cells.Item(16, 13).Value = -169383.33
$ws.Cells.Item(16, 14).Value = -6266.6665
# Row 75
$ws.Cells.Item(75, 8).Value = 16250
$ws.Cells.Item(75, 10).Value = 16250
$ws.Cells.Item(75, 12).Value = 16250
$ws.Cells.Item(75, 14).Value = -18246
# Row 78
$ws.Cells.Item(78, 8).Value = 16250
$ws.Cells.Item(78, 10).Value = 16250
$ws.Cells.Item(78, 12).Value = 48750
$ws.Cells.Item(78, 14).Value = -58734
# Row 86
$ws.Cells.Item(86, 8).Value = 21531.871
$ws.Cells.Item(86, 9).Value = 39024.438
$ws.Cells.Item(86, 10).Value = 2873.1333
$ws.Cells.Item(86, 11).Value = 39024.438
$ws.Cells.Item(86, 12).Value = 2873.1333
$ws.Cells.Item(86, 13).Value = -37901.438
$ws.Cells.Item(86, 14).Value = -5119.1333
# Row 88
$ws.Cells.Item(88, 8).Value = 18612.312
$ws.Cells.Item(88, 9).Value = 15000
$ws.Cells.Item(88, 10).Value = 18853.133
$ws.Cells.Item(88, 11).Value = 15000
$ws.Cells.Item(88, 12).Value = 18853.133
$ws.Cells.Item(88, 13).Value = -14594
$ws.Cells.Item(88, 14).Value = -19665.133
# Row 89
$ws.Cells.Item(89, 8).Value = 21531.871
$ws.Cells.Item(89, 9).Value = 39024.438
$ws.Cells.Item(89, 10).Value = 2873.1333
$ws.Cells.Item(89, 11).Value = 195122.19
$ws.Cells.Item(89, 12).Value = 14365.6665
$ws.Cells.Item(89, 13).Value = -189506.19
$ws.Cells.Item(89, 14).Value = -25597.6665
# Row 91
$ws.Cells.Item(91, 8).Value = 18612.312
$ws.Cells.Item(91, 9).Value = 15000
$ws.Cells.Item(91, 10).Value = 18853.133
$ws.Cells.Item(91, 11).Value = 15000
$ws.Cells.Item(91, 12).Value = 18853.133
$ws.Cells.Item(91, 13).Value = -13596
$ws.Cells.Item(91, 14).Value = -21661.133
# Row 113
$ws.Cells.Item(113, 8).Value = 60351.89
$ws.Cells.Item(113, 9).Value = 169670.33
$ws.Cells.Item(113, 10).Value = 5692.6665
$ws.Cells.Item(113, 11).Value = 169670.33
$ws.Cells.Item(113, 12).Value = 5692.6665
$ws.Cells.Item(113, 13).Value = -167500.33
$ws.Cells.Item(113, 14).Value = -10032.6665

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Cells.Item(8, 8).Value = 154.5625
$ws.Cells.Item(8, 9).Value = 154.5625
$ws.Cells.Item(8, 11).Value = 463.6875
$ws.Cells.Item(8, 13).Value = -324.6875
# Row 64
$ws.Cells.Item(64, 8).Value = 2804.348
$ws.Cells.Item(64, 9).Value = 3940
$ws.Cells.Item(64, 10).Value = 2488.889
$ws.Cells.Item(64, 11).Value = 11820
$ws.Cells.Item(64, 12).Value = 7466.667
$ws.Cells.Item(64, 13).Value = -11550
$ws.Cells.Item(64, 14).Value = -8006.667
# Row 67
$ws.Cells.Item(67, 8).Value = 2804.348
$ws.Cells.Item(67, 9).Value = 3940
$ws.Cells.Item(67, 10).Value = 2488.889
$ws.Cells.Item(67, 11).Value = 11820
$ws.Cells.Item(67, 12).Value = 7466.667
$ws.Cells.Item(67, 13).Value = -10884
$ws.Cells.Item(67, 14).Value = -9338.667000000001
# Row 80
$ws.Cells.Item(80, 8).Value = 2392.6
$ws.Cells.Item(80, 9).Value = 950
$ws.Cells.Item(80, 10).Value = 2614.5386
$ws.Cells.Item(80, 11).Value = 2850
$ws.Cells.Item(80, 12).Value = 7843.6158
$ws.Cells.Item(80, 13).Value = -1914
$ws.Cells.Item(80, 14).Value = -9715.6158
# Row 83
$ws.Cells.Item(83, 8).Value = 2392.6
$ws.Cells.Item(83, 9).Value = 950
$ws.Cells.Item(83, 10).Value = 2614.5386
$ws.Cells.Item(83, 11).Value = 8550
$ws.Cells.Item(83, 12).Value = 23530.8474
$ws.Cells.Item(83, 13).Value = -3870
$ws.Cells.Item(83, 14).Value = -32890.8474
# Row 108
$ws.Cells.Item(108, 8).Value = 940.5714
$ws.Cells.Item(108, 9).Value = 940.5714
$ws.Cells.Item(108, 11).Value = 2821.7142
$ws.Cells.Item(108, 13).Value = 58.28579999999965

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 39
$ws.Cells.Item(39, 8).Value = 20550
$ws.Cells.Item(39, 10).Value = 20550
$ws.Cells.Item(39, 12).Value = 20550
$ws.Cells.Item(39, 14).Value = -21614
# Row 69
$ws.Cells.Item(69, 8).Value = 29000
$ws.Cells.Item(69, 10).Value = 29000
$ws.Cells.Item(69, 12).Value = 29000
$ws.Cells.Item(69, 14).Value = -30498
# Row 72
$ws.Cells.Item(72, 8).Value = 29000
$ws.Cells.Item(72, 10).Value = 29000
$ws.Cells.Item(72, 12).Value = 87000
$ws.Cells.Item(72, 14).Value = -94488
# Row 92
$ws.Cells.Item(92, 8).Value = 8384
$ws.Cells.Item(92, 10).Value = 8384
$ws.Cells.Item(92, 12).Value = 8384
$ws.Cells.Item(92, 14).Value = -12128
# Row 101
$ws.Cells.Item(101, 8).Value = 55885.668
$ws.Cells.Item(101, 10).Value = 55885.668
$ws.Cells.Item(101, 12).Value = 55885.668
$ws.Cells.Item(101, 14).Value = -62375.668
# Row 107
$ws.Cells.Item(107, 8).Value = 1807
$ws.Cells.Item(107, 9).Value = 2641.25
$ws.Cells.Item(107, 11).Value = 2641.25
$ws.Cells.Item(107, 13).Value = -721.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 43
$ws.Cells.Item(43, 8).Value = 10000
$ws.Cells.Item(43, 9).Value = 2000
$ws.Cells.Item(43, 10).Value = 11600
$ws.Cells.Item(43, 11).Value = 2000
$ws.Cells.Item(43, 12).Value = 11600
$ws.Cells.Item(43, 13).Value = -1807
$ws.Cells.Item(43, 14).Value = -11986
# Row 55
$ws.Cells.Item(55, 8).Value = 148.85715
$ws.Cells.Item(55, 10).Value = 183.46666
$ws.Cells.Item(55, 12).Value = 183.46666
$ws.Cells.Item(55, 14).Value = -529.46666
# Row 100
$ws.Cells.Item(100, 8).Value = 2725.3333
$ws.Cells.Item(100, 9).Value = 2670.4
$ws.Cells.Item(100, 10).Value = 3000
$ws.Cells.Item(100, 11).Value = 2670.4
$ws.Cells.Item(100, 12).Value = 3000
$ws.Cells.Item(100, 13).Value = -2129.4
$ws.Cells.Item(100, 14).Value = -4082
# Row 122
$ws.Cells.Item(122, 8).Value = 8820.789000000001
$ws.Cells.Item(122, 9).Value = 8820.789000000001
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 26462.367
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -24012.367
$ws.Cells.Item(122, 14).ClearContents() | Out-Null

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 75
$ws.Cells.Item(75, 8).Value = 33695.555
$ws.Cells.Item(75, 10).Value = 33695.555
$ws.Cells.Item(75, 12).Value = 33695.555
$ws.Cells.Item(75, 14).Value = -35567.555
# Row 78
$ws.Cells.Item(78, 8).Value = 33695.555
$ws.Cells.Item(78, 10).Value = 33695.555
$ws.Cells.Item(78, 12).Value = 101086.665
$ws.Cells.Item(78, 14).Value = -110446.665
# Row 95
$ws.Cells.Item(95, 8).Value = 31580.5
$ws.Cells.Item(95, 10).Value = 31580.5
$ws.Cells.Item(95, 12).Value = 31580.5
$ws.Cells.Item(95, 14).Value = -37072.5
# Row 104
$ws.Cells.Item(104, 8).Value = 10370
$ws.Cells.Item(104, 10).Value = 10370
$ws.Cells.Item(104, 12).Value = 10370
$ws.Cells.Item(104, 14).Value = -17358
# Row 136
$ws.Cells.Item(136, 8).Value = 4921.0293
$ws.Cells.Item(136, 9).Value = 1306.6957
$ws.Cells.Item(136, 10).Value = 12478.272
$ws.Cells.Item(136, 11).Value = 3920.0871
$ws.Cells.Item(136, 12).Value = 37434.81600000001
$ws.Cells.Item(136, 13).Value = -1370.0871
$ws.Cells.Item(136, 14).Value = -42534.81600000001

Write-Output "done"